$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new blank rows before the old row 81 ("button.close"), which will
# push it (and everything after) down to row 90, making room for the new
# "scheduleRun.*" field rows 81-89.
$ws.Range("A81:A89").EntireRow.Insert()

# Row 81: scheduleRun.search
$ws.Cells.Item(81, 1).Value = "scheduleRun.search"
$ws.Cells.Item(81, 2).Value = "Nhập email hoặc tên bản mẫu để tìm kiếm"
$ws.Cells.Item(81, 3).Value = "Input email or template name to search"
$ws.Rows.Item(81).RowHeight = 15

# Row 82: scheduleRun.id
$ws.Cells.Item(82, 1).Value = "scheduleRun.id"
$ws.Cells.Item(82, 2).Value = "Mã"
$ws.Cells.Item(82, 3).Value = "ID"
$ws.Rows.Item(82).RowHeight = 15

# Row 83: scheduleRun.email
$ws.Cells.Item(83, 1).Value = "scheduleRun.email"
$ws.Cells.Item(83, 2).Value = "Email"
$ws.Cells.Item(83, 3).Value = "Email"
$ws.Rows.Item(83).RowHeight = 15

# Row 84: scheduleRun.emailTos
$ws.Cells.Item(84, 1).Value = "scheduleRun.emailTos"
$ws.Cells.Item(84, 2).Value = "Gửi tới"
$ws.Cells.Item(84, 3).Value = "To"
$ws.Rows.Item(84).RowHeight = 15

# Row 85: scheduleRun.proxy
$ws.Cells.Item(85, 1).Value = "scheduleRun.proxy"
$ws.Cells.Item(85, 2).Value = "Proxy"
$ws.Cells.Item(85, 3).Value = "Proxy"
$ws.Rows.Item(85).RowHeight = 15

# Row 86: scheduleRun.schedule
$ws.Cells.Item(86, 1).Value = "scheduleRun.schedule"
$ws.Cells.Item(86, 2).Value = "Lịch trình"
$ws.Cells.Item(86, 3).Value = "Schedule"
$ws.Rows.Item(86).RowHeight = 15

# Row 87: scheduleRun.template
$ws.Cells.Item(87, 1).Value = "scheduleRun.template"
$ws.Cells.Item(87, 2).Value = "Bản mẫu"
$ws.Cells.Item(87, 3).Value = "Template"
$ws.Rows.Item(87).RowHeight = 15

# Row 88: scheduleRun.enable
$ws.Cells.Item(88, 1).Value = "scheduleRun.enable"
$ws.Cells.Item(88, 2).Value = "Trạng thái"
$ws.Cells.Item(88, 3).Value = "Enable"
$ws.Rows.Item(88).RowHeight = 15

# Row 89: scheduleRun.actions (uses a distinct Calibri 12pt font, slightly
# taller row to match the author's edit)
$ws.Cells.Item(89, 1).Value = "scheduleRun.actions"
$ws.Cells.Item(89, 2).Value = "Hành động"
$ws.Cells.Item(89, 3).Value = "Actions"
$ws.Cells.Item(89, 1).Font.Name = "Calibri"
$ws.Cells.Item(89, 1).Font.Size = 12
$ws.Rows.Item(89).RowHeight = 15.65

# Update the selection to match the author's final cursor position.
$ws.Range("B90").Select()
